$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(3, 1).Value = "Globo"
$ws.Cells.Item(3, 2).Value = "Inter TV Rural"
$ws.Cells.Item(3, 3).Value = "Agricultura"
$ws.Cells.Item(3, 4).Value = "2025-03-31T19:36"
$ws.Cells.Item(3, 5).Value = "Negativo"
$ws.Cells.Item(3, 6).Value = "Sem Nota"
$ws.Cells.Item(3, 7).Value = "teste"
